$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 634
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 902
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 902
$ws.Range("M9").Value = -331
$ws.Range("N9").Value = -1240
$ws.Range("H19").Value = 38967.816
$ws.Range("J19").Value = 53279.625
$ws.Range("L19").Value = 53279.625
$ws.Range("N19").Value = -53629.625
$ws.Range("H46").Value = 1842.3334
$ws.Range("I46").Value = 1695.6666
$ws.Range("K46").Value = 5086.9998
$ws.Range("M46").Value = -4967.9998
$ws.Range("H60").Value = 1842.3334
$ws.Range("I60").Value = 1695.6666
$ws.Range("K60").Value = 5086.9998
$ws.Range("M60").Value = -4602.9998
$ws.Range("H61").Value = 4302.8
$ws.Range("I61").Value = 4128.5
$ws.Range("K61").Value = 12385.5
$ws.Range("M61").Value = -12213.5
$ws.Range("H92").Value = 3142.0715
$ws.Range("I92").Value = 3142.0715
$ws.Range("K92").Value = 3142.0715
$ws.Range("M92").Value = -1894.0715
$ws.Range("H107").Value = 1007.5333
$ws.Range("I107").Value = 1007.5333
$ws.Range("K107").Value = 1007.5333
$ws.Range("M107").Value = 912.4666999999999
$ws.Range("H116").Value = 6011.9165
$ws.Range("I116").Value = 4655.5713
$ws.Range("K116").Value = 4655.5713
$ws.Range("M116").Value = -1213.5713
$ws.Range("H132").Value = 13340.8
$ws.Range("I132").Value = 14277.421
$ws.Range("J132").Value = 11723
$ws.Range("K132").Value = 42832.263
$ws.Range("L132").Value = 35169
$ws.Range("M132").Value = -40302.263
$ws.Range("N132").Value = -40229
$ws.Range("H137").Value = 11480.435
$ws.Range("I137").Value = 2454.3125
$ws.Range("J137").Value = 32111.572
$ws.Range("K137").Value = 7362.9375
$ws.Range("L137").Value = 96334.716
$ws.Range("M137").Value = -4812.9375
$ws.Range("N137").Value = -101434.716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("H30").Value = 5000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("H32").Value = 3908753.8
$ws.Range("I32").Value = 4809082.5
$ws.Range("K32").Value = 4809082.5
$ws.Range("M32").Value = -4808795.5
$ws.Range("H61").Value = 1026108.1
$ws.Range("I61").Value = 3672.8
$ws.Range("K61").Value = 3672.8
$ws.Range("M61").Value = -3460.8
$ws.Range("H97").Value = 736.3
$ws.Range("I97").Value = 613.2941
$ws.Range("J97").Value = 1433.3334
$ws.Range("K97").Value = 613.2941
$ws.Range("L97").Value = 1433.3334
$ws.Range("M97").Value = -117.2941
$ws.Range("N97").Value = -2425.3334
$ws.Range("H102").Value = 11268.5
$ws.Range("I102").Value = 7260
$ws.Range("J102").Value = 15277
$ws.Range("K102").Value = 7260
$ws.Range("L102").Value = 15277
$ws.Range("M102").Value = -5638
$ws.Range("N102").Value = -18521
$ws.Range("H103").Value = 41500
$ws.Range("J103").Value = 41500
$ws.Range("L103").Value = 41500
$ws.Range("N103").Value = -43844
$ws.Range("H110").Value = 8541.857
$ws.Range("I110").Value = 9632.166999999999
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 9632.166999999999
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -7587.166999999999
$ws.Range("N110").Value = -6090
$ws.Range("H132").Value = 3487364.8
$ws.Range("I132").Value = 1679.5454
$ws.Range("J132").Value = 8280182
$ws.Range("K132").Value = 5038.6362
$ws.Range("L132").Value = 24840546
$ws.Range("M132").Value = -2508.6362
$ws.Range("N132").Value = -24845606
$ws.Range("H136").Value = 1026108.1
$ws.Range("I136").Value = 3672.8
$ws.Range("K136").Value = 11018.4
$ws.Range("M136").Value = -8468.400000000001
$ws.Range("M3").ClearContents()
$ws.Range("N30").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1375
$ws.Range("I11").Value = 354
$ws.Range("J11").Value = 2651.25
$ws.Range("K11").Value = 354
$ws.Range("L11").Value = 2651.25
$ws.Range("M11").Value = -214
$ws.Range("N11").Value = -2931.25
$ws.Range("H12").Value = 660.8
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 813.5
$ws.Range("K12").Value = 50
$ws.Range("L12").Value = 813.5
$ws.Range("M12").Value = 118
$ws.Range("N12").Value = -1149.5
$ws.Range("H22").Value = 6333.3335
$ws.Range("I22").Value = 6333.3335
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 6333.3335
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -6160.3335
$ws.Range("H24").Value = 700
$ws.Range("I24").Value = 700
$ws.Range("K24").Value = 700
$ws.Range("M24").Value = -465
$ws.Range("H25").Value = 14000
$ws.Range("I25").Value = 14000
$ws.Range("K25").Value = 14000
$ws.Range("M25").Value = -13765
$ws.Range("H34").Value = 2026
$ws.Range("J34").Value = 2026
$ws.Range("L34").Value = 2026
$ws.Range("N34").Value = -2254
$ws.Range("H37").Value = 2205.2
$ws.Range("I37").Value = 1008.6667
$ws.Range("J37").Value = 4000
$ws.Range("K37").Value = 1008.6667
$ws.Range("L37").Value = 4000
$ws.Range("M37").Value = -871.6667
$ws.Range("N37").Value = -4274
$ws.Range("H99").Value = 13517.85
$ws.Range("I99").Value = 14974
$ws.Range("J99").Value = 5266.3335
$ws.Range("K99").Value = 14974
$ws.Range("L99").Value = 5266.3335
$ws.Range("M99").Value = -13476
$ws.Range("N99").Value = -8262.333500000001
$ws.Range("H134").Value = 61052.9
$ws.Range("I134").Value = 63376
$ws.Range("K134").Value = 190128
$ws.Range("M134").Value = -187593
$ws.Range("N22").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 784
$ws.Range("I22").Value = 410.5
$ws.Range("K22").Value = 410.5
$ws.Range("M22").Value = -60.5
$ws.Range("H31").Value = 9405.777
$ws.Range("I31").Value = 925.9643
$ws.Range("J31").Value = 39085.125
$ws.Range("K31").Value = 925.9643
$ws.Range("L31").Value = 39085.125
$ws.Range("M31").Value = -630.9643
$ws.Range("N31").Value = -39675.125
$ws.Range("H34").Value = 9405.777
$ws.Range("I34").Value = 925.9643
$ws.Range("J34").Value = 39085.125
$ws.Range("K34").Value = 925.9643
$ws.Range("L34").Value = 39085.125
$ws.Range("M34").Value = -723.9643
$ws.Range("N34").Value = -39489.125
$ws.Range("H132").Value = 21906424
$ws.Range("I132").Value = 2023.9
$ws.Range("K132").Value = 6071.700000000001
$ws.Range("M132").Value = -3541.700000000001
$ws.Range("H134").Value = 18523492
$ws.Range("I134").Value = 3050.919
$ws.Range("J134").Value = 58832690
$ws.Range("K134").Value = 9152.757
$ws.Range("L134").Value = 176498070
$ws.Range("M134").Value = -6617.757
$ws.Range("N134").Value = -176503140

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1792.3334
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("H117").Value = 356.7143
$ws.Range("I117").Value = 216.16667
$ws.Range("K117").Value = 648.50001
$ws.Range("M117").Value = 2793.49999
$ws.Range("M114").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -7
$ws.Range("H21").Value = 18333.334
$ws.Range("H29").Value = 15000
$ws.Range("J29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15580
$ws.Range("H30").Value = 18333.334
$ws.Range("H41").Value = 3017
$ws.Range("I41").Value = 3017
$ws.Range("K41").Value = 3017
$ws.Range("M41").Value = -2662
$ws.Range("H55").Value = 29000
$ws.Range("J55").Value = 29000
$ws.Range("L55").Value = 29000
$ws.Range("N55").Value = -29654
$ws.Range("H99").Value = 8687.5
$ws.Range("I99").Value = 8687.5
$ws.Range("K99").Value = 8687.5
$ws.Range("M99").Value = -6441.5
$ws.Range("H132").Value = 552754.8
$ws.Range("I132").Value = 2970.5925
$ws.Range("K132").Value = 8911.7775
$ws.Range("M132").Value = -6381.7775
$ws.Range("N18").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 13421.866
$ws.Range("I68").Value = 10722.059
$ws.Range("J68").Value = 16952.385
$ws.Range("K68").Value = 10722.059
$ws.Range("L68").Value = 16952.385
$ws.Range("M68").Value = -9973.058999999999
$ws.Range("N68").Value = -18450.385
$ws.Range("H71").Value = 13421.866
$ws.Range("I71").Value = 10722.059
$ws.Range("J71").Value = 16952.385
$ws.Range("K71").Value = 53610.295
$ws.Range("L71").Value = 84761.92499999999
$ws.Range("M71").Value = -49866.295
$ws.Range("N71").Value = -92249.92499999999
$ws.Range("H122").Value = 7142.7144
$ws.Range("I122").Value = 6649.75
$ws.Range("K122").Value = 19949.25
$ws.Range("M122").Value = -17499.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678
$ws.Range("H122").Value = 2745.3
$ws.Range("I122").Value = 2587.5715
$ws.Range("K122").Value = 7762.7145
$ws.Range("M122").Value = -5312.7145
